$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

$ws.Cells.Item($row, 4).Value = Get-Date -Year 2022 -Month 7 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112037
$ws.Cells.Item($row, 7).Value = "Cebollín"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 8000
$ws.Cells.Item($row, 12).Value = 8000
$ws.Cells.Item($row, 13).Value = 8000
$ws.Cells.Item($row, 14).Value = "$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 2667
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
